$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 4 data rows (rows 10-13); this also updates the sheet
# dimension and shared-string usage automatically.
$ws.Rows("10:13").Delete()

# Row 2: ECs -> FAPs (Sema6d-Tyrobp)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema6d"
$ws.Range("C2").Value = "Tyrobp"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 58.62893533333334
$ws.Range("H2").Value = 175.886806
$ws.Range("I2").Value = 0.5702456571409142
$ws.Range("J2").Value = 0.5702456571409142
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.8937520000000001
$ws.Range("N2").Value = 2.681256
$ws.Range("O2").Value = 0.001982631228033859
$ws.Range("P2").Value = 0.001982631228033858
$ws.Range("Q2").Value = 52.39972821203735
$ws.Range("R2").Value = 471.5975539083361
$ws.Range("S2").Value = 0.001130586847498265
$ws.Range("T2").Value = 0.001130586847498265

# Row 3: ECs -> Resolving-Mac
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema6d"
$ws.Range("C3").Value = "Tyrobp"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 58.62893533333334
$ws.Range("H3").Value = 175.886806
$ws.Range("I3").Value = 0.5702456571409142
$ws.Range("J3").Value = 0.5702456571409142
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 449.8970896666667
$ws.Range("N3").Value = 1349.691269
$ws.Range("O3").Value = 0.9980173687719661
$ws.Range("P3").Value = 0.9980173687719661
$ws.Range("Q3").Value = 26376.98737672187
$ws.Range("R3").Value = 237392.8863904969
$ws.Range("S3").Value = 0.569115070293416
$ws.Range("T3").Value = 0.569115070293416

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sema6d"
$ws.Range("C4").Value = "Tyrobp"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 14.28901333333333
$ws.Range("H4").Value = 42.86704
$ws.Range("I4").Value = 0.1389799721218762
$ws.Range("J4").Value = 0.1389799721218763
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.8937520000000001
$ws.Range("N4").Value = 2.681256
$ws.Range("O4").Value = 0.001982631228033859
$ws.Range("P4").Value = 0.001982631228033858
$ws.Range("Q4").Value = 12.77083424469334
$ws.Range("R4").Value = 114.93750820224
$ws.Range("S4").Value = 0.0002755460328001069
$ws.Range("T4").Value = 0.0002755460328001069

# Row 5: FAPs -> Resolving-Mac
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema6d"
$ws.Range("C5").Value = "Tyrobp"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.28901333333333
$ws.Range("H5").Value = 42.86704
$ws.Range("I5").Value = 0.1389799721218762
$ws.Range("J5").Value = 0.1389799721218763
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 449.8970896666667
$ws.Range("N5").Value = 1349.691269
$ws.Range("O5").Value = 0.9980173687719661
$ws.Range("P5").Value = 0.9980173687719661
$ws.Range("Q5").Value = 6428.585512874864
$ws.Range("R5").Value = 57857.26961587377
$ws.Range("S5").Value = 0.1387044260890761
$ws.Range("T5").Value = 0.1387044260890762

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Sema6d"
$ws.Range("C6").Value = "Tyrobp"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 27.27518533333334
$ws.Range("H6").Value = 81.82555600000001
$ws.Range("I6").Value = 0.265288050953297
$ws.Range("J6").Value = 0.2652880509532971
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.8937520000000001
$ws.Range("N6").Value = 2.681256
$ws.Range("O6").Value = 0.001982631228033859
$ws.Range("P6").Value = 0.001982631228033858
$ws.Range("Q6").Value = 24.37725144203734
$ws.Range("R6").Value = 219.395262978336
$ws.Range("S6").Value = 0.0005259683742442442
$ws.Range("T6").Value = 0.0005259683742442441

# Row 7: MuSCs -> Resolving-Mac
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Sema6d"
$ws.Range("C7").Value = "Tyrobp"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 27.27518533333334
$ws.Range("H7").Value = 81.82555600000001
$ws.Range("I7").Value = 0.265288050953297
$ws.Range("J7").Value = 0.2652880509532971
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 449.8970896666667
$ws.Range("N7").Value = 1349.691269
$ws.Range("O7").Value = 0.9980173687719661
$ws.Range("P7").Value = 0.9980173687719661
$ws.Range("Q7").Value = 12271.02650158562
$ws.Range("R7").Value = 110439.2385142706
$ws.Range("S7").Value = 0.2647620825790528
$ws.Range("T7").Value = 0.2647620825790529

# Row 8: Resolving-Mac -> FAPs
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Sema6d"
$ws.Range("C8").Value = "Tyrobp"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.620337
$ws.Range("H8").Value = 7.861011
$ws.Range("I8").Value = 0.02548631978391236
$ws.Range("J8").Value = 0.02548631978391236
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.8937520000000001
$ws.Range("N8").Value = 2.681256
$ws.Range("O8").Value = 0.001982631228033859
$ws.Range("P8").Value = 0.001982631228033858
$ws.Range("Q8").Value = 2.341931434424
$ws.Range("R8").Value = 21.077382909816
$ws.Range("S8").Value = 0.00005052997349124178
$ws.Range("T8").Value = 0.00005052997349124178

# Row 9: Resolving-Mac -> Resolving-Mac
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Sema6d"
$ws.Range("C9").Value = "Tyrobp"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.620337
$ws.Range("H9").Value = 7.861011
$ws.Range("I9").Value = 0.02548631978391236
$ws.Range("J9").Value = 0.02548631978391236
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 449.8970896666667
$ws.Range("N9").Value = 1349.691269
$ws.Range("O9").Value = 0.9980173687719661
$ws.Range("P9").Value = 0.9980173687719661
$ws.Range("Q9").Value = 1178.881990245884
$ws.Range("R9").Value = 10609.93791221296
$ws.Range("S9").Value = 0.02543578981042111
$ws.Range("T9").Value = 0.02543578981042112
